$wb = $excel.ActiveWorkbook

# --- Build Romania sheet from a copy of the UK sheet (which has the full
#     repeater item superset), then trim to the Romania-specific subset. ---
$uk = $wb.Worksheets.Item("UK")
$uk.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$romania = $wb.Worksheets.Item($wb.Worksheets.Count)
$romania.Name = "Romania"

# Drop the FC32AR / FC32DR rows (rows 10-11) that UK has but Romania doesn't.
$romania.Range("A10:A11").EntireRow.Delete() | Out-Null

# --- Build Slovakia sheet the same way. ---
$uk2 = $wb.Worksheets.Item("UK")
$uk2.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

$slovakia.Range("A10:A11").EntireRow.Delete() | Out-Null

# Market names for both new sheets.
$romania.Range("B2").Value = "Romania Market"
$slovakia.Range("B2").Value = "Slovakia Market"

# User story / input value for both new sheets.
$romania.Range("B4").Value = "NGC-4307/T3536/T3543"
$slovakia.Range("B4").Value = "NGC-4306/T3562/T3575"

# Romania is not the selected/active sheet; restore the default (non-full-row)
# selection state that a freshly duplicated sheet would have.
$romania.Range("B4").Select() | Out-Null
$slovakia.Range("B6").Select() | Out-Null

# --- UK sheet: selection moves off B4 onto a full-row selection (A1:XFD1048576),
#     and it is no longer the tab-selected sheet once Slovakia becomes active. ---
$uk.Range("A1:XFD1048576").Select() | Out-Null

# --- Slovakia becomes the active (selected) tab. ---
$slovakia.Activate() | Out-Null
